$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24: set date, clear "Anderson" note in column E
$ws.Range("A24").Value = 41838
$ws.Range("E24").Value = ""

# Row 26: "Realizado ajustes..." must be introduced first so it becomes
# shared string index 45 (matching original authoring order)
$ws.Range("A26").Value = 41840
$ws.Range("C26").Value = "Realizado ajustes no ""CE"" da calculadora"
$ws.Range("B26").Value = "Implementado uma nova interface e testes"
$ws.Range("E26").Value = "Germán"

# Row 25: fill in remaining standup notes
$ws.Range("A25").Value = 41839
$ws.Range("C25").Value = "Implementado uma nova interface e testes"
$ws.Range("D25").Value = "Problemas no ""CE"" para número"
$ws.Range("E25").Value = "Germán"

# Row 27: set date, clear "Ronaldo" note in column E
$ws.Range("A27").Value = 41841
$ws.Range("E27").Value = ""

# Row 28: clear remaining placeholder values
$ws.Range("A28").Value = ""
$ws.Range("E28").Value = ""

# Column C width adjustment (auto-fit to the new, wider content)
$ws.Columns.Item(3).ColumnWidth = 63

# Selection moves to B28
$ws.Range("B28").Select()
